# Update country data figures on the Jordan Summary sheet.
# The affected cells hold numbers stored as text in the workbook, so a
# leading apostrophe is used to force a text ("quote prefix") entry just
# like typing the value directly into Excel would, keeping the cell's
# original General number format instead of converting it to a real
# numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people)
$ws.Range("B11").Value = "'23.19"
$ws.Range("C11").Value = "'2.05"
$ws.Range("D11").Value = "'25.25"

# Employment (% of total)
$ws.Range("B12").Value = "'40.97"
$ws.Range("C12").Value = "'24.71"
$ws.Range("D12").Value = "'65.68"

# Enterprises (% of total)
$ws.Range("B14").Value = "'91.47"
$ws.Range("D14").Value = "'99.57"
